# edit.ps1 - applies the "end" style slide fixes described by the diff:
#   - the big background picture ("Рисунок 6") is renamed to "Рисунок 7"
#     and brought forward in z-order (from the very back to just behind
#     the TITLE text box, i.e. in front of the caption / logo / QR code)
#   - several shapes receive a tiny re-measure / size correction
#     (the classic +0.14pt / +1800 EMU growth that shows up after the
#     deck is round-tripped through PowerPoint)
#   - the TITLE placeholder text is recolored from black to white and
#     shrunk slightly (3400 -> 3300), fixing the invisible-text-on-dark-
#     background bug, and its box grows a touch to fit

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Grab stable references to every shape on the slide up front (by their
# current, pre-edit z-order position) so index shuffling caused by the
# later re-ordering step doesn't affect the rest of the script.
$bigPic   = $s.Shapes.Item(1)   # "Рисунок 6" - big background picture
$rect5    = $s.Shapes.Item(2)   # "Прямоугольник 5" - footer caption
$ellipse  = $s.Shapes.Item(3)   # "Изображение 40" - round logo
$qrPic    = $s.Shapes.Item(4)   # "Рисунок 5" - QR code picture
$titleBox = $s.Shapes.Item(5)   # "" - TITLE text box

# --- Rename the background picture ("Рисунок 6" -> "Рисунок 7") -----------
$bigPic.Name = "Рисунок 7"

# --- Tiny size corrections (EMU -> points, 12700 EMU per point) -----------
$bigPic.Width   = 7565400 / 12700
$bigPic.Height  = 5141160 / 12700

$rect5.Width    = 2997000 / 12700
$rect5.Height   = 267480 / 12700

$ellipse.Width  = 254880 / 12700
$ellipse.Height = 254880 / 12700

$qrPic.Width    = 1287360 / 12700
$qrPic.Height   = 1131480 / 12700

$titleBox.Width  = 4744800 / 12700
$titleBox.Height = 1023120 / 12700

# --- Fix the TITLE text formatting (black -> white, 34pt -> 33pt) ---------
$titleRange = $titleBox.TextFrame.TextRange
$titleRange.Font.Size = 33
$titleRange.Font.Color.RGB = 0xFFFFFF

# --- Bring the background picture forward so it now sits just behind -----
# --- the TITLE box instead of being the very backmost shape ---------------
$bigPic.ZOrder(2)   # msoBringForward
$bigPic.ZOrder(2)
$bigPic.ZOrder(2)
